$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G:G").Delete()
